# "anik separated Asser regisyter" - split the Asset Code column (B) so the
# single category letter (O/C/F/T/M/I/S/E/P ...) is replaced with a fixed
# "7-4" segment for the first batch of assets (rows 2-151) and "7-n" for the
# rest (rows 152-301): FRC-HQ-SLM-<X>-<YY>-<NNNN> -> FRC-SLM-7-4-<YY>-<NNNN>
# or FRC-SLM-7-n-<YY>-<NNNN>.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 301
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $old = [string]$cell.Value2

    if ($old -match '^FRC-HQ-SLM-[A-Za-z]-(\d{2}-\d{4})$') {
        $tail = $matches[1]
        if ($row -le 151) {
            $cell.Value = "FRC-SLM-7-4-$tail"
        } else {
            $cell.Value = "FRC-SLM-7-n-$tail"
        }
    }
}

# Two numeric cells that moved very slightly (floating point noise on the
# "Cost of Assets Sold" / "Current Balance" pair, still summing to the
# original "Price") also changed in the same commit.
$ws.Range("Q18").Value = 30145.79439252336
$ws.Range("R18").Value = -0.004392523358546896

$ws.Range("Q55").Value = 13789.79591836735
$ws.Range("R55").Value = 0.004081632649103994
